$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "xltablediff.py --key ID test1old.xlsx test1new.xlsx --out test1diff.xlsx"

$ws.Range("B5").Value = "Diff test:"
$ws.Range("C5").Value = "xltablediff.py  --key ID test1old.xlsx test1new.xlsx --out test1diff.xlsx"

$ws.Range("B6").Value = "Merge test:"
$ws.Range("C6").Value = "xltablediff.py  --key ID --merge Color test1old.xlsx test1new.xlsx --out test1merge.xlsx"

$ws.Range("B7").Value = "Append test:"
$ws.Range("C7").Value = "xltablediff.py  --key ID --append test1old.xlsx test1new.xlsx --out test1append.xlsx"
